$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")
$ws.Range("C1").Formula = "=LEN(A1)"
$ws.Range("C2:C65").Formula = "=LEN(A2)"
$ws.Range("C66:C129").Formula = "=LEN(A66)"
$ws.Range("C130:C183").Formula = "=LEN(A130)"
[void]$ws.Range("B7").Select()
